# Auto-generated Excel COM-interop script applying the Behemoth_Profits value updates
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 28 (hunk 0)
$ws.Cells.Item(28, 8).Value = 1367.375  # H28: was 1367.5
$ws.Cells.Item(28, 9).Value = 981  # I28: was 981.3333
$ws.Cells.Item(28, 11).Value = 981  # K28: was 981.3333
$ws.Cells.Item(28, 13).Value = -496  # M28: was -496.3333
# row 127 (hunk 1)
$ws.Cells.Item(127, 8).Value = 10913.571  # H127: was 7843.8
$ws.Cells.Item(127, 9).Value = 697.5  # I127: was 687.6
$ws.Cells.Item(127, 11).Value = 2092.5  # K127: was 2062.8
$ws.Cells.Item(127, 13).Value = 2867.5  # M127: was 2897.2
# row 132 (hunk 2)
$ws.Cells.Item(132, 8).Value = 696.1786  # H132: was 731.39624
$ws.Cells.Item(132, 9).Value = 698.46295  # I132: was 735.1961
$ws.Cells.Item(132, 11).Value = 2095.38885  # K132: was 2205.5883
$ws.Cells.Item(132, 13).Value = 434.6111500000002  # M132: was 324.4117000000001
# row 135 (hunk 3)
$ws.Cells.Item(135, 8).Value = 2638.8  # H135: was 1772.875
$ws.Cells.Item(135, 9).Value = 3199.75  # I135: was 1772.875
$ws.Cells.Item(135, 10).Value = 395  # J135: was 0
$ws.Cells.Item(135, 11).Value = 28797.75  # K135: was 15955.875
$ws.Cells.Item(135, 12).Value = 3555  # L135: was 0
$ws.Cells.Item(135, 13).Value = -26262.75  # M135: was -13420.875
$ws.Cells.Item(135, 14).Value = -8625  # N135: add
# row 138 (hunk 4)
$ws.Cells.Item(138, 8).Value = 2392  # H138: was 2400.04
$ws.Cells.Item(138, 10).Value = 2448.5269  # J138: was 2456.4788
$ws.Cells.Item(138, 12).Value = 7345.5807  # L138: was 7369.4364
$ws.Cells.Item(138, 14).Value = -17625.5807  # N138: was -17649.4364

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 2 (hunk 5)
$ws.Cells.Item(2, 8).Value = 1491.875  # H2: was 1519.4286
$ws.Cells.Item(2, 9).Value = 1491.875  # I2: was 1519.4286
$ws.Cells.Item(2, 11).Value = 1491.875  # K2: was 1519.4286
$ws.Cells.Item(2, 13).Value = -1378.875  # M2: was -1406.4286
# row 4 (hunk 6)
$ws.Cells.Item(4, 8).Value = 1394.6842  # H4: was 1921
$ws.Cells.Item(4, 9).Value = 1479.8572  # I4: was 1835.909
$ws.Cells.Item(4, 10).Value = 1156.2  # J4: was 2389
$ws.Cells.Item(4, 11).Value = 1479.8572  # K4: was 1835.909
$ws.Cells.Item(4, 12).Value = 1156.2  # L4: was 2389
$ws.Cells.Item(4, 13).Value = -1363.8572  # M4: was -1719.909
$ws.Cells.Item(4, 14).Value = -1388.2  # N4: was -2621
# row 28 (hunk 7)
$ws.Cells.Item(28, 8).Value = 55498.8  # H28: was 61499.75
$ws.Cells.Item(28, 9).Value = 41623.5  # I28: was 44999.668
$ws.Cells.Item(28, 11).Value = 41623.5  # K28: was 44999.668
$ws.Cells.Item(28, 13).Value = -41431.5  # M28: was -44807.668
# row 32 (hunk 8)
$ws.Cells.Item(32, 8).Value = 6352736.5  # H32: was 5838795.5
$ws.Cells.Item(32, 9).Value = 8642458  # I32: was 8642462
$ws.Cells.Item(32, 10).Value = 28743.047  # J32: was 31200.143
$ws.Cells.Item(32, 11).Value = 8642458  # K32: was 8642462
$ws.Cells.Item(32, 12).Value = 28743.047  # L32: was 31200.143
$ws.Cells.Item(32, 13).Value = -8642171  # M32: was -8642175
$ws.Cells.Item(32, 14).Value = -29317.047  # N32: was -31774.143
# row 61 (hunk 9)
$ws.Cells.Item(61, 8).Value = 30003326  # H61: was 32612236
$ws.Cells.Item(61, 9).Value = 25002716  # I61: was 26318664
$ws.Cells.Item(61, 10).Value = 50005764  # J61: was 62506704
$ws.Cells.Item(61, 11).Value = 25002716  # K61: was 26318664
$ws.Cells.Item(61, 12).Value = 50005764  # L61: was 62506704
$ws.Cells.Item(61, 13).Value = -25002504  # M61: was -26318452
$ws.Cells.Item(61, 14).Value = -50006188  # N61: was -62507128
# row 99 (hunk 10)
$ws.Cells.Item(99, 8).Value = 55498.8  # H99: was 61499.75
$ws.Cells.Item(99, 9).Value = 41623.5  # I99: was 44999.668
$ws.Cells.Item(99, 11).Value = 41623.5  # K99: was 44999.668
$ws.Cells.Item(99, 13).Value = -38628.5  # M99: was -42004.668
# row 116 (hunk 11)
$ws.Cells.Item(116, 8).Value = 1491.875  # H116: was 1519.4286
$ws.Cells.Item(116, 9).Value = 1491.875  # I116: was 1519.4286
$ws.Cells.Item(116, 11).Value = 1491.875  # K116: was 1519.4286
$ws.Cells.Item(116, 13).Value = 802.125  # M116: was 774.5714
# row 122 (hunk 12)
$ws.Cells.Item(122, 8).Value = 3743.7942  # H122: was 3871.5312
$ws.Cells.Item(122, 9).Value = 3448.4546  # I122: was 3530.182
$ws.Cells.Item(122, 10).Value = 4285.25  # J122: was 4622.5
$ws.Cells.Item(122, 11).Value = 10345.3638  # K122: was 10590.546
$ws.Cells.Item(122, 12).Value = 12855.75  # L122: was 13867.5
$ws.Cells.Item(122, 13).Value = -7895.363799999999  # M122: was -8140.545999999998
$ws.Cells.Item(122, 14).Value = -17755.75  # N122: was -18767.5
# row 132 (hunk 13)
$ws.Cells.Item(132, 8).Value = 7579854.5  # H132: was 8134412
$ws.Cells.Item(132, 9).Value = 11496847  # I132: was 12823302
$ws.Cells.Item(132, 11).Value = 34490541  # K132: was 38469906
$ws.Cells.Item(132, 13).Value = -34488011  # M132: was -38467376
# row 136 (hunk 14)
$ws.Cells.Item(136, 8).Value = 30003326  # H136: was 32612236
$ws.Cells.Item(136, 9).Value = 25002716  # I136: was 26318664
$ws.Cells.Item(136, 10).Value = 50005764  # J136: was 62506704
$ws.Cells.Item(136, 11).Value = 75008148  # K136: was 78955992
$ws.Cells.Item(136, 12).Value = 150017292  # L136: was 187520112
$ws.Cells.Item(136, 13).Value = -75005598  # M136: was -78953442
$ws.Cells.Item(136, 14).Value = -150022392  # N136: was -187525212

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 3 (hunk 15)
$ws.Cells.Item(3, 8).Value = 1491.875  # H3: was 1519.4286
$ws.Cells.Item(3, 9).Value = 1491.875  # I3: was 1519.4286
$ws.Cells.Item(3, 11).Value = 1491.875  # K3: was 1519.4286
$ws.Cells.Item(3, 13).Value = -1377.875  # M3: was -1405.4286
# row 82 (hunk 16)
$ws.Cells.Item(82, 8).Value = 17387.9  # H82: was 21653.545
$ws.Cells.Item(82, 9).Value = 4863  # I82: was 3742.625
$ws.Cells.Item(82, 10).Value = 46612.668  # J82: was 69416
$ws.Cells.Item(82, 11).Value = 4863  # K82: was 3742.625
$ws.Cells.Item(82, 12).Value = 46612.668  # L82: was 69416
$ws.Cells.Item(82, 13).Value = -4480  # M82: was -3359.625
$ws.Cells.Item(82, 14).Value = -47378.668  # N82: was -70182
# row 85 (hunk 17)
$ws.Cells.Item(85, 8).Value = 17387.9  # H85: was 21653.545
$ws.Cells.Item(85, 9).Value = 4863  # I85: was 3742.625
$ws.Cells.Item(85, 10).Value = 46612.668  # J85: was 69416
$ws.Cells.Item(85, 11).Value = 4863  # K85: was 3742.625
$ws.Cells.Item(85, 12).Value = 46612.668  # L85: was 69416
$ws.Cells.Item(85, 13).Value = -3537  # M85: was -2416.625
$ws.Cells.Item(85, 14).Value = -49264.668  # N85: was -72068
# row 86 (hunk 18)
$ws.Cells.Item(86, 8).Value = 2390.88  # H86: was 2612.318
$ws.Cells.Item(86, 9).Value = 2240.0557  # I86: was 2435.3125
$ws.Cells.Item(86, 10).Value = 2778.7144  # J86: was 3084.3333
$ws.Cells.Item(86, 11).Value = 2240.0557  # K86: was 2435.3125
$ws.Cells.Item(86, 12).Value = 2778.7144  # L86: was 3084.3333
$ws.Cells.Item(86, 13).Value = -1117.0557  # M86: was -1312.3125
$ws.Cells.Item(86, 14).Value = -5024.7144  # N86: was -5330.3333
# row 89 (hunk 19)
$ws.Cells.Item(89, 8).Value = 2390.88  # H89: was 2612.318
$ws.Cells.Item(89, 9).Value = 2240.0557  # I89: was 2435.3125
$ws.Cells.Item(89, 10).Value = 2778.7144  # J89: was 3084.3333
$ws.Cells.Item(89, 11).Value = 11200.2785  # K89: was 12176.5625
$ws.Cells.Item(89, 12).Value = 13893.572  # L89: was 15421.6665
$ws.Cells.Item(89, 13).Value = -5584.2785  # M89: was -6560.5625
$ws.Cells.Item(89, 14).Value = -25125.572  # N89: was -26653.6665
# row 105 (hunk 20)
$ws.Cells.Item(105, 8).Value = 2362.3845  # H105: was 2575.0833
$ws.Cells.Item(105, 9).Value = 1566.6666  # I105: was 1918
$ws.Cells.Item(105, 11).Value = 1566.6666  # K105: was 1918
$ws.Cells.Item(105, 13).Value = 180.3334  # M105: was -171

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 2 (hunk 21)
$ws.Cells.Item(2, 8).Value = 30466.666  # H2: was 36680
$ws.Cells.Item(2, 10).Value = 0  # J2: was 46000
$ws.Cells.Item(2, 12).Value = 0  # L2: was 46000
$ws.Cells.Item(2, 14).ClearContents()  # N2: was -46226
# row 31 (hunk 22)
$ws.Cells.Item(31, 8).Value = 1375537.2  # H31: was 1063241.9
$ws.Cells.Item(31, 9).Value = 2115.7778  # I31: was 1852.1428
$ws.Cells.Item(31, 10).Value = 2920636.5  # J31: was 2920674
$ws.Cells.Item(31, 11).Value = 2115.7778  # K31: was 1852.1428
$ws.Cells.Item(31, 12).Value = 2920636.5  # L31: was 2920674
$ws.Cells.Item(31, 13).Value = -1820.7778  # M31: was -1557.1428
$ws.Cells.Item(31, 14).Value = -2921226.5  # N31: was -2921264
# row 34 (hunk 23)
$ws.Cells.Item(34, 8).Value = 1375537.2  # H34: was 1063241.9
$ws.Cells.Item(34, 9).Value = 2115.7778  # I34: was 1852.1428
$ws.Cells.Item(34, 10).Value = 2920636.5  # J34: was 2920674
$ws.Cells.Item(34, 11).Value = 2115.7778  # K34: was 1852.1428
$ws.Cells.Item(34, 12).Value = 2920636.5  # L34: was 2920674
$ws.Cells.Item(34, 13).Value = -1913.7778  # M34: was -1650.1428
$ws.Cells.Item(34, 14).Value = -2921040.5  # N34: was -2921078
# row 132 (hunk 24)
$ws.Cells.Item(132, 8).Value = 4706.6216  # H132: was 4824.0557
$ws.Cells.Item(132, 9).Value = 1423.5  # I132: was 1440.84
$ws.Cells.Item(132, 10).Value = 12466.728  # J132: was 12513.182
$ws.Cells.Item(132, 11).Value = 4270.5  # K132: was 4322.52
$ws.Cells.Item(132, 12).Value = 37400.18399999999  # L132: was 37539.546
$ws.Cells.Item(132, 13).Value = -1740.5  # M132: was -1792.52
$ws.Cells.Item(132, 14).Value = -42460.18399999999  # N132: was -42599.546

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 87 (hunk 25)
$ws.Cells.Item(87, 8).Value = 12916.667  # H87: was 17250
$ws.Cells.Item(87, 9).Value = 4375  # I87: was 10875
$ws.Cells.Item(87, 11).Value = 13125  # K87: was 32625
$ws.Cells.Item(87, 13).Value = -11877  # M87: was -31377
# row 88 (hunk 26)
$ws.Cells.Item(88, 8).Value = 4534.5  # H88: was 4570.3076
$ws.Cells.Item(88, 10).Value = 4850  # J88: was 4863.636
$ws.Cells.Item(88, 12).Value = 14550  # L88: was 14590.908
$ws.Cells.Item(88, 14).Value = -15406  # N88: was -15446.908
# row 90 (hunk 27)
$ws.Cells.Item(90, 8).Value = 12916.667  # H90: was 17250
$ws.Cells.Item(90, 9).Value = 4375  # I90: was 10875
$ws.Cells.Item(90, 11).Value = 39375  # K90: was 97875
$ws.Cells.Item(90, 13).Value = -33135  # M90: was -91635
# row 91 (hunk 28)
$ws.Cells.Item(91, 8).Value = 4534.5  # H91: was 4570.3076
$ws.Cells.Item(91, 10).Value = 4850  # J91: was 4863.636
$ws.Cells.Item(91, 12).Value = 14550  # L91: was 14590.908
$ws.Cells.Item(91, 14).Value = -17514  # N91: was -17554.908
# row 94 (hunk 29)
$ws.Cells.Item(94, 8).Value = 3483.3333  # H94: was 2947.913
$ws.Cells.Item(94, 9).Value = 900  # I94: was 828.8570999999999
$ws.Cells.Item(94, 10).Value = 4000  # J94: was 3875
$ws.Cells.Item(94, 11).Value = 2700  # K94: was 2486.5713
$ws.Cells.Item(94, 12).Value = 12000  # L94: was 11625
$ws.Cells.Item(94, 13).Value = -2024  # M94: was -1810.5713
$ws.Cells.Item(94, 14).Value = -13352  # N94: was -12977
# row 113 (hunk 30)
$ws.Cells.Item(113, 8).Value = 1806.9333  # H113: was 1854.7142
$ws.Cells.Item(113, 10).Value = 1801.8572  # J113: was 1852.9231
$ws.Cells.Item(113, 12).Value = 5405.571599999999  # L113: was 5558.7693
$ws.Cells.Item(113, 14).Value = -9745.571599999999  # N113: was -9898.7693

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 18 (hunk 31)
$ws.Cells.Item(18, 8).Value = 42500  # H18: was 55000
$ws.Cells.Item(18, 9).Value = 30000  # I18: was 0
$ws.Cells.Item(18, 11).Value = 30000  # K18: was 0
$ws.Cells.Item(18, 13).Value = -29707  # M18: add
# row 113 (hunk 32)
$ws.Cells.Item(113, 8).Value = 4368.45  # H113: was 4456.5264
$ws.Cells.Item(113, 9).Value = 3663.3  # I113: was 3770.889
$ws.Cells.Item(113, 11).Value = 3663.3  # K113: was 3770.889
$ws.Cells.Item(113, 13).Value = -1493.3  # M113: was -1600.889
# row 132 (hunk 33)
$ws.Cells.Item(132, 8).Value = 58827280  # H132: was 66670704
$ws.Cells.Item(132, 9).Value = 90913610  # I132: was 100004870
$ws.Cells.Item(132, 10).Value = 2335.3333  # J132: was 2362.4
$ws.Cells.Item(132, 11).Value = 272740830  # K132: was 300014610
$ws.Cells.Item(132, 12).Value = 7005.999899999999  # L132: was 7087.200000000001
$ws.Cells.Item(132, 13).Value = -272738300  # M132: was -300012080
$ws.Cells.Item(132, 14).Value = -12065.9999  # N132: was -12147.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 16 (hunk 34)
$ws.Cells.Item(16, 8).Value = 878.36  # H16: was 854.37036
$ws.Cells.Item(16, 9).Value = 771  # I16: was 752.9583
$ws.Cells.Item(16, 11).Value = 771  # K16: was 752.9583
$ws.Cells.Item(16, 13).Value = -601  # M16: was -582.9583
# row 40 (hunk 35)
$ws.Cells.Item(40, 8).Value = 4468.407  # H40: was 4298.241
$ws.Cells.Item(40, 9).Value = 4185.1665  # I40: was 4070.158
$ws.Cells.Item(40, 10).Value = 5034.8887  # J40: was 4731.6
$ws.Cells.Item(40, 11).Value = 4185.1665  # K40: was 4070.158
$ws.Cells.Item(40, 12).Value = 5034.8887  # L40: was 4731.6
$ws.Cells.Item(40, 13).Value = -4049.1665  # M40: was -3934.158
$ws.Cells.Item(40, 14).Value = -5306.8887  # N40: was -5003.6
# row 46 (hunk 36)
$ws.Cells.Item(46, 8).Value = 3352.238  # H46: was 3462.35
$ws.Cells.Item(46, 9).Value = 2459.077  # I46: was 2568.1667
$ws.Cells.Item(46, 11).Value = 2459.077  # K46: was 2568.1667
$ws.Cells.Item(46, 13).Value = -2271.077  # M46: was -2380.1667
# row 56 (hunk 37)
$ws.Cells.Item(56, 8).Value = 22498.334  # H56: was 29999.5
$ws.Cells.Item(56, 9).Value = 11500  # I56: was 20000
$ws.Cells.Item(56, 10).Value = 44495  # J56: was 39999
$ws.Cells.Item(56, 11).Value = 11500  # K56: was 20000
$ws.Cells.Item(56, 12).Value = 44495  # L56: was 39999
$ws.Cells.Item(56, 13).Value = -10809  # M56: was -19309
$ws.Cells.Item(56, 14).Value = -45877  # N56: was -41381
# row 61 (hunk 38)
$ws.Cells.Item(61, 8).Value = 2321.5293  # H61: was 2275.45
$ws.Cells.Item(61, 9).Value = 1831.1333  # I61: was 1861.6666
$ws.Cells.Item(61, 11).Value = 1831.1333  # K61: was 1861.6666
$ws.Cells.Item(61, 13).Value = -1629.1333  # M61: was -1659.6666
# row 74 (hunk 39)
$ws.Cells.Item(74, 8).Value = 16064600  # H74: was 13395333
$ws.Cells.Item(74, 9).Value = 26699332  # I74: was 20036750
$ws.Cells.Item(74, 11).Value = 26699332  # K74: was 20036750
$ws.Cells.Item(74, 13).Value = -26698334  # M74: was -20035752
# row 77 (hunk 40)
$ws.Cells.Item(77, 8).Value = 16064600  # H77: was 13395333
$ws.Cells.Item(77, 9).Value = 26699332  # I77: was 20036750
$ws.Cells.Item(77, 11).Value = 80097996  # K77: was 60110250
$ws.Cells.Item(77, 13).Value = -80093004  # M77: was -60105258
# row 82 (hunk 41)
$ws.Cells.Item(82, 8).Value = 2390.5454  # H82: was 1782.5
$ws.Cells.Item(82, 9).Value = 1532.6666  # I82: was 1397.25
$ws.Cells.Item(82, 10).Value = 3420  # J82: was 2002.6428
$ws.Cells.Item(82, 11).Value = 1532.6666  # K82: was 1397.25
$ws.Cells.Item(82, 12).Value = 3420  # L82: was 2002.6428
$ws.Cells.Item(82, 13).Value = -1171.6666  # M82: was -1036.25
$ws.Cells.Item(82, 14).Value = -4142  # N82: was -2724.6428
# row 85 (hunk 42)
$ws.Cells.Item(85, 8).Value = 2390.5454  # H85: was 1782.5
$ws.Cells.Item(85, 9).Value = 1532.6666  # I85: was 1397.25
$ws.Cells.Item(85, 10).Value = 3420  # J85: was 2002.6428
$ws.Cells.Item(85, 11).Value = 1532.6666  # K85: was 1397.25
$ws.Cells.Item(85, 12).Value = 3420  # L85: was 2002.6428
$ws.Cells.Item(85, 13).Value = -284.6666  # M85: was -149.25
$ws.Cells.Item(85, 14).Value = -5916  # N85: was -4498.6428
# row 113 (hunk 43)
$ws.Cells.Item(113, 8).Value = 2321.5293  # H113: was 2275.45
$ws.Cells.Item(113, 9).Value = 1831.1333  # I113: was 1861.6666
$ws.Cells.Item(113, 11).Value = 1831.1333  # K113: was 1861.6666
$ws.Cells.Item(113, 13).Value = 338.8667  # M113: was 308.3334

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 26 (hunk 44)
$ws.Cells.Item(26, 8).Value = 0  # H26: was 16499.5
$ws.Cells.Item(26, 9).Value = 0  # I26: was 16499.5
$ws.Cells.Item(26, 11).Value = 0  # K26: was 16499.5
$ws.Cells.Item(26, 13).ClearContents()  # M26: was -16206.5
# row 122 (hunk 45)
$ws.Cells.Item(122, 8).Value = 3843.1428  # H122: was 3501.5
$ws.Cells.Item(122, 9).Value = 3843.1428  # I122: was 3501.5
$ws.Cells.Item(122, 11).Value = 11529.4284  # K122: was 10504.5
$ws.Cells.Item(122, 13).Value = -9079.428400000001  # M122: was -8054.5
# row 126 (hunk 46)
$ws.Cells.Item(126, 8).Value = 1340.8928  # H126: was 1311.862
$ws.Cells.Item(126, 9).Value = 1128.3043  # I126: was 1102.0834
$ws.Cells.Item(126, 11).Value = 3384.9129  # K126: was 3306.2502
$ws.Cells.Item(126, 13).Value = -914.9129000000003  # M126: was -836.2501999999999
# row 132 (hunk 47)
$ws.Cells.Item(132, 8).Value = 3092621.2  # H132: was 3340073.2
$ws.Cells.Item(132, 9).Value = 5836.9  # I132: was 6544.3887
$ws.Cells.Item(132, 11).Value = 17510.7  # K132: was 19633.1661
$ws.Cells.Item(132, 13).Value = -14980.7  # M132: was -17103.1661
